$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.01222466666667
$ws.Range("H2").Value = 135.036674
$ws.Range("I2").Value = 0.7482903203664146
$ws.Range("J2").Value = 0.7482903203664146
$ws.Range("M2").Value = 0.6692693333333334
$ws.Range("N2").Value = 2.007808
$ws.Range("O2").Value = 0.004126561180566838
$ws.Range("P2").Value = 0.004126561180566839
$ws.Range("Q2").Value = 30.12530159451023
$ws.Range("R2").Value = 271.127714350592
$ws.Range("S2").Value = 0.00308786578781797
$ws.Range("T2").Value = 0.00308786578781797
$ws.Range("G3").Value = 45.01222466666667
$ws.Range("H3").Value = 135.036674
$ws.Range("I3").Value = 0.7482903203664146
$ws.Range("J3").Value = 0.7482903203664146
$ws.Range("O3").Value = 0.9916964991825307
$ws.Range("P3").Value = 0.9916964991825309
$ws.Range("Q3").Value = 7239.722088402417
$ws.Range("R3").Value = 65157.49879562177
$ws.Range("S3").Value = 0.7420768910795477
$ws.Range("T3").Value = 0.7420768910795479
$ws.Range("G4").Value = 45.01222466666667
$ws.Range("H4").Value = 135.036674
$ws.Range("I4").Value = 0.7482903203664146
$ws.Range("J4").Value = 0.7482903203664146
$ws.Range("M4").Value = 0.5637343333333333
$ws.Range("N4").Value = 1.691203
$ws.Range("O4").Value = 0.003475856580040611
$ws.Range("P4").Value = 0.003475856580040611
$ws.Range("Q4").Value = 25.37493646431356
$ws.Range("R4").Value = 228.374428178822
$ws.Range("S4").Value = 0.002600949833826299
$ws.Range("T4").Value = 0.002600949833826299
$ws.Range("G5").Value = 45.01222466666667
$ws.Range("H5").Value = 135.036674
$ws.Range("I5").Value = 0.7482903203664146
$ws.Range("J5").Value = 0.7482903203664146
$ws.Range("M5").Value = 0.1137056666666667
$ws.Range("N5").Value = 0.341117
$ws.Range("O5").Value = 0.0007010830568617209
$ws.Range("P5").Value = 0.0007010830568617211
$ws.Range("Q5").Value = 5.118145013873111
$ws.Range("R5").Value = 46.063305124858
$ws.Range("S5").Value = 0.0005246136652225224
$ws.Range("T5").Value = 0.0005246136652225225
$ws.Range("G6").Value = 2.766295666666667
$ws.Range("H6").Value = 8.298887000000001
$ws.Range("I6").Value = 0.04598733535094824
$ws.Range("J6").Value = 0.04598733535094825
$ws.Range("M6").Value = 0.6692693333333334
$ws.Range("N6").Value = 2.007808
$ws.Range("O6").Value = 0.004126561180566838
$ws.Range("P6").Value = 0.004126561180566839
$ws.Range("Q6").Value = 1.851396856632889
$ws.Range("R6").Value = 16.662571709696
$ws.Range("S6").Value = 0.0001897695528569321
$ws.Range("T6").Value = 0.0001897695528569321
$ws.Range("G7").Value = 2.766295666666667
$ws.Range("H7").Value = 8.298887000000001
$ws.Range("I7").Value = 0.04598733535094824
$ws.Range("J7").Value = 0.04598733535094825
$ws.Range("O7").Value = 0.9916964991825307
$ws.Range("P7").Value = 0.9916964991825309
$ws.Range("Q7").Value = 444.9282831348147
$ws.Range("R7").Value = 4004.354548213333
$ws.Range("S7").Value = 0.04560547947426841
$ws.Range("T7").Value = 0.04560547947426842
$ws.Range("G8").Value = 2.766295666666667
$ws.Range("H8").Value = 8.298887000000001
$ws.Range("I8").Value = 0.04598733535094824
$ws.Range("J8").Value = 0.04598733535094825
$ws.Range("M8").Value = 0.5637343333333333
$ws.Range("N8").Value = 1.691203
$ws.Range("O8").Value = 0.003475856580040611
$ws.Range("P8").Value = 0.003475856580040611
$ws.Range("Q8").Value = 1.559455843451222
$ws.Range("R8").Value = 14.035102591061
$ws.Range("S8").Value = 0.0001598453821781276
$ws.Range("T8").Value = 0.0001598453821781277
$ws.Range("G9").Value = 2.766295666666667
$ws.Range("H9").Value = 8.298887000000001
$ws.Range("I9").Value = 0.04598733535094824
$ws.Range("J9").Value = 0.04598733535094825
$ws.Range("M9").Value = 0.1137056666666667
$ws.Range("N9").Value = 0.341117
$ws.Range("O9").Value = 0.0007010830568617209
$ws.Range("P9").Value = 0.0007010830568617211
$ws.Range("Q9").Value = 0.3145434929754444
$ws.Range("R9").Value = 2.830891436779
$ws.Range("S9").Value = 0.00003224094164476787
$ws.Range("T9").Value = 0.00003224094164476789
$ws.Range("G10").Value = 12.37490333333333
$ws.Range("H10").Value = 37.12471
$ws.Range("I10").Value = 0.2057223442826371
$ws.Range("J10").Value = 0.2057223442826371
$ws.Range("M10").Value = 0.6692693333333334
$ws.Range("N10").Value = 2.007808
$ws.Range("O10").Value = 0.004126561180566838
$ws.Range("P10").Value = 0.004126561180566839
$ws.Range("Q10").Value = 8.282143303964446
$ws.Range("R10").Value = 74.53928973568001
$ws.Range("S10").Value = 0.0008489258398919366
$ws.Range("T10").Value = 0.0008489258398919369
$ws.Range("G11").Value = 12.37490333333333
$ws.Range("H11").Value = 37.12471
$ws.Range("I11").Value = 0.2057223442826371
$ws.Range("J11").Value = 0.2057223442826371
$ws.Range("O11").Value = 0.9916964991825307
$ws.Range("P11").Value = 0.9916964991825309
$ws.Range("Q11").Value = 1990.367320603099
$ws.Range("R11").Value = 17913.30588542789
$ws.Range("S11").Value = 0.2040141286287145
$ws.Range("T11").Value = 0.2040141286287146
$ws.Range("G12").Value = 12.37490333333333
$ws.Range("H12").Value = 37.12471
$ws.Range("I12").Value = 0.2057223442826371
$ws.Range("J12").Value = 0.2057223442826371
$ws.Range("M12").Value = 0.5637343333333333
$ws.Range("N12").Value = 1.691203
$ws.Range("O12").Value = 0.003475856580040611
$ws.Range("P12").Value = 0.003475856580040611
$ws.Range("Q12").Value = 6.976157880681112
$ws.Range("R12").Value = 62.78542092613
$ws.Range("S12").Value = 0.0007150613640361842
$ws.Range("T12").Value = 0.0007150613640361844
$ws.Range("G13").Value = 12.37490333333333
$ws.Range("H13").Value = 37.12471
$ws.Range("I13").Value = 0.2057223442826371
$ws.Range("J13").Value = 0.2057223442826371
$ws.Range("M13").Value = 0.1137056666666667
$ws.Range("N13").Value = 0.341117
$ws.Range("O13").Value = 0.0007010830568617209
$ws.Range("P13").Value = 0.0007010830568617211
$ws.Range("Q13").Value = 1.407096633452222
$ws.Range("R13").Value = 12.66386970107
$ws.Range("S13").Value = 0.0001442284499944306
$ws.Range("T13").Value = 0.0001442284499944306
